# Apply cash-flow figure updates on the "GLW" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GLW")

# Row 4 - Non Cash Items (Other)
$ws.Range("B4").Value = -598000000.0
$ws.Range("F4").Value = 453000000.0

# Row 6 - Change in inventories
$ws.Range("B6").Value = 534000000.0
$ws.Range("C6").Value = 423000000.0
$ws.Range("D6").Value = 280000000.0
$ws.Range("E6").Value = 12000000.0
$ws.Range("F6").Value = -206000000.0

# Row 8 - Change in payables and accrued liability
$ws.Range("B8").Value = 4731000000.0
$ws.Range("C8").Value = 4709000000.0
$ws.Range("D8").Value = 5122000000.0

# Row 12 - Capital expenditures
$ws.Range("B12").Value = -1084000000.0

# Row 13 - Net Aquisitions (was blank, now populated)
$ws.Range("B13").Value = -20000000.0

# Row 19 - Equity Repurchase (Common, Net)
$ws.Range("B19").Value = 168000000.0
$ws.Range("F19").Value = -723000000.0

# Row 30 - Assets Liabilities Change (Total)
$ws.Range("B30").Value = 546000000.0
$ws.Range("F30").Value = -60000000.0

# Row 32 - Issuance/Purchase of Shares
$ws.Range("B32").Value = 168000000.0
$ws.Range("F32").Value = -723000000.0
